# This edit removes two checklist items ("Fill the Skill Assessment Survey
# and commit it on the GitHub" and "Pre-implementation Questionnaire") and
# relocates Word's auto "last edit" bookmark (_GoBack) so that it now sits
# at the start of what becomes the new first checklist item
# ("Understand the application ...").

$d = $word.ActiveDocument

# --- Step 1: relocate the hidden "_GoBack" bookmark -----------------------
# It currently sits just before "Post-implementation Questionnaire"; after
# the edit it should sit just before "Understand the application " (the
# paragraph that becomes the first bullet once "Fill the ... GitHub" is
# removed).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $d.Paragraphs.Item(4)          # "Understand the application ..."
$insertPoint = $d.Range($target.Range.Start, $target.Range.Start)
$d.Bookmarks.Add("_GoBack", $insertPoint)

# --- Step 2: delete the "Fill the ... GitHub" bullet item -----------------
# Deleting the paragraph's Range (which includes its trailing paragraph
# mark) merges it away entirely, leaving the following paragraph intact.
$d.Paragraphs.Item(3).Range.Delete()

# --- Step 3: delete the "Pre-implementation Questionnaire" bullet item ----
# After step 2, this item has shifted up to become paragraph 4.
$d.Paragraphs.Item(4).Range.Delete()
